# Auto-generated edit script applying scheduled market-data refresh to Seraph_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 254.53847
$ws.Range("J33").Value = 287.5
$ws.Range("L33").Value = 287.5
$ws.Range("N33").Value = -745.5
$ws.Range("H40").Value = 2333.5557
$ws.Range("I40").Value = 1999.909
$ws.Range("J40").Value = 2857.8572
$ws.Range("K40").Value = 1999.909
$ws.Range("L40").Value = 2857.8572
$ws.Range("M40").Value = -1824.909
$ws.Range("N40").Value = -3207.8572
$ws.Range("H62").Value = 1104.5
$ws.Range("I62").Value = 1104.5
$ws.Range("K62").Value = 1104.5
$ws.Range("M62").Value = -480.5
$ws.Range("H65").Value = 1104.5
$ws.Range("I65").Value = 1104.5
$ws.Range("K65").Value = 5522.5
$ws.Range("M65").Value = -2402.5
$ws.Range("H70").Value = 54789.69
$ws.Range("I70").Value = 1999.5
$ws.Range("J70").Value = 64387.91
$ws.Range("K70").Value = 5998.5
$ws.Range("L70").Value = 193163.73
$ws.Range("M70").Value = -5728.5
$ws.Range("N70").Value = -193703.73
$ws.Range("H73").Value = 54789.69
$ws.Range("I73").Value = 1999.5
$ws.Range("J73").Value = 64387.91
$ws.Range("K73").Value = 5998.5
$ws.Range("L73").Value = 193163.73
$ws.Range("M73").Value = -5062.5
$ws.Range("N73").Value = -195035.73

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 2859762.5
$ws.Range("J6").Value = 10001500
$ws.Range("L6").Value = 10001500
$ws.Range("N6").Value = -10001846
$ws.Range("H13").Value = 7501448.5
$ws.Range("I13").Value = 10001265
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 10001265
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = -10001121
$ws.Range("N13").Value = -2288
$ws.Range("H32").Value = 4209.3486
$ws.Range("I32").Value = 2535.7646
$ws.Range("J32").Value = 10531.777
$ws.Range("K32").Value = 2535.7646
$ws.Range("L32").Value = 10531.777
$ws.Range("M32").Value = -2248.7646
$ws.Range("N32").Value = -11105.777
$ws.Range("H45").Value = 2899.5
$ws.Range("J45").Value = 4200
$ws.Range("L45").Value = 4200
$ws.Range("N45").Value = -4954
$ws.Range("H61").Value = 9557.333000000001
$ws.Range("I61").Value = 10127
$ws.Range("K61").Value = 10127
$ws.Range("M61").Value = -9915
$ws.Range("H110").Value = 6947298.5
$ws.Range("I110").Value = 27778594
$ws.Range("J110").Value = 3533.3333
$ws.Range("K110").Value = 27778594
$ws.Range("L110").Value = 3533.3333
$ws.Range("M110").Value = -27776549
$ws.Range("N110").Value = -7623.3333
$ws.Range("H132").Value = 6840.2
$ws.Range("I132").Value = 7650.25
$ws.Range("K132").Value = 22950.75
$ws.Range("M132").Value = -20420.75
$ws.Range("H136").Value = 9557.333000000001
$ws.Range("I136").Value = 10127
$ws.Range("K136").Value = 30381
$ws.Range("M136").Value = -27831

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 16668556
$ws.Range("I7").Value = 25000334
$ws.Range("K7").Value = 25000334
$ws.Range("M7").Value = -25000221
$ws.Range("H105").Value = 4171303.8
$ws.Range("I105").Value = 7580243.5
$ws.Range("J105").Value = 4821.6665
$ws.Range("K105").Value = 7580243.5
$ws.Range("L105").Value = 4821.6665
$ws.Range("M105").Value = -7578496.5
$ws.Range("N105").Value = -8315.666499999999
$ws.Range("H134").Value = 1918.7
$ws.Range("I134").Value = 2199.125
$ws.Range("J134").Value = 797
$ws.Range("K134").Value = 6597.375
$ws.Range("L134").Value = 2391
$ws.Range("M134").Value = -4062.375
$ws.Range("N134").Value = -7461

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 2000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = -2340
$ws.Range("H24").Value = 2000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = ""
$ws.Range("N24").Value = -2340
$ws.Range("H31").Value = 3854.8823
$ws.Range("I31").Value = 2780.0435
$ws.Range("K31").Value = 2780.0435
$ws.Range("M31").Value = -2485.0435
$ws.Range("H34").Value = 3854.8823
$ws.Range("I34").Value = 2780.0435
$ws.Range("K34").Value = 2780.0435
$ws.Range("M34").Value = -2578.0435
$ws.Range("H58").Value = 2893.7693
$ws.Range("I58").Value = 1786
$ws.Range("J58").Value = 5386.25
$ws.Range("K58").Value = 1786
$ws.Range("L58").Value = 5386.25
$ws.Range("M58").Value = -1583
$ws.Range("N58").Value = -5792.25
$ws.Range("H62").Value = 203000.5
$ws.Range("I62").Value = 6002
$ws.Range("K62").Value = 6002
$ws.Range("M62").Value = -5378
$ws.Range("H65").Value = 203000.5
$ws.Range("I65").Value = 6002
$ws.Range("K65").Value = 30010
$ws.Range("M65").Value = -26890
$ws.Range("H99").Value = 10508.5625
$ws.Range("I99").Value = 6371.778
$ws.Range("K99").Value = 6371.778
$ws.Range("M99").Value = -4873.778
$ws.Range("H126").Value = 10508.5625
$ws.Range("I126").Value = 6371.778
$ws.Range("K126").Value = 19115.334
$ws.Range("M126").Value = -16645.334
$ws.Range("H134").Value = 3996.8333
$ws.Range("I134").Value = 3995.3333
$ws.Range("K134").Value = 11985.9999
$ws.Range("M134").Value = -9450.999899999999
$ws.Range("H136").Value = 2893.7693
$ws.Range("I136").Value = 1786
$ws.Range("J136").Value = 5386.25
$ws.Range("K136").Value = 5358
$ws.Range("L136").Value = 16158.75
$ws.Range("M136").Value = -2808
$ws.Range("N136").Value = -21258.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 38661700
$ws.Range("I4").Value = 48747132
$ws.Range("K4").Value = 146241396
$ws.Range("M4").Value = -146241284
$ws.Range("H13").Value = 20
$ws.Range("I13").Value = 20
$ws.Range("K13").Value = 60
$ws.Range("M13").Value = 108

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 20000
$ws.Range("I62").Value = 20000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 20000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -19314
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 20000
$ws.Range("I65").Value = 20000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 60000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -56568
$ws.Range("N65").Value = ""
$ws.Range("H122").Value = 113384.555
$ws.Range("I122").Value = 1743.6666
$ws.Range("K122").Value = 5230.9998
$ws.Range("M122").Value = -2780.9998
$ws.Range("H123").Value = 40000.5
$ws.Range("J123").Value = 40000.5
$ws.Range("L123").Value = 40000.5
$ws.Range("N123").Value = -44900.5
$ws.Range("H132").Value = 3149.6
$ws.Range("I132").Value = 2944
$ws.Range("K132").Value = 8832
$ws.Range("M132").Value = -6302

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4510.5
$ws.Range("I7").Value = 3930
$ws.Range("J7").Value = 5555.4
$ws.Range("K7").Value = 3930
$ws.Range("L7").Value = 5555.4
$ws.Range("M7").Value = -3818
$ws.Range("N7").Value = -5779.4
$ws.Range("H12").Value = 2000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = -2340
$ws.Range("H63").Value = 20000
$ws.Range("I63").Value = 20000
$ws.Range("K63").Value = 20000
$ws.Range("M63").Value = -19251
$ws.Range("H66").Value = 20000
$ws.Range("I66").Value = 20000
$ws.Range("K66").Value = 60000
$ws.Range("M66").Value = -56256
$ws.Range("H122").Value = 4004
$ws.Range("I122").Value = 4004
$ws.Range("K122").Value = 12012
$ws.Range("M122").Value = -9562
$ws.Range("H126").Value = 4510.5
$ws.Range("I126").Value = 3930
$ws.Range("J126").Value = 5555.4
$ws.Range("K126").Value = 11790
$ws.Range("L126").Value = 16666.2
$ws.Range("M126").Value = -9320
$ws.Range("N126").Value = -21606.2
$ws.Range("H136").Value = 7696.1665
$ws.Range("J136").Value = 8042.25
$ws.Range("L136").Value = 24126.75
$ws.Range("N136").Value = -29226.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 17796.8
$ws.Range("I11").Value = 17000
$ws.Range("J11").Value = 17996
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 17996
$ws.Range("M11").Value = -16858
$ws.Range("N11").Value = -18280
$ws.Range("H15").Value = 15000
$ws.Range("J15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("N15").Value = -15576
$ws.Range("H62").Value = 4725.2666
$ws.Range("I62").Value = 3972
$ws.Range("K62").Value = 3972
$ws.Range("M62").Value = -3348
$ws.Range("H65").Value = 4725.2666
$ws.Range("I65").Value = 3972
$ws.Range("K65").Value = 19860
$ws.Range("M65").Value = -16740
$ws.Range("H112").Value = 16000
$ws.Range("J112").Value = 16000
$ws.Range("L112").Value = 16000
$ws.Range("N112").Value = -18954
$ws.Range("H113").Value = 1810.4445
$ws.Range("J113").Value = 3249
$ws.Range("L113").Value = 9747
$ws.Range("N113").Value = -14087
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("H136").Value = 2545.8572
$ws.Range("I136").Value = 1386.625
$ws.Range("J136").Value = 4091.5
$ws.Range("K136").Value = 4159.875
$ws.Range("L136").Value = 12274.5
$ws.Range("M136").Value = -1609.875
$ws.Range("N136").Value = -17374.5
